$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the current row 94 (before the old row 95),
# which shifts the old rows 95-119 down to 97-121 and extends the used range to R121.
$ws.Rows("95:96").Insert()

# New row 95: Camote, "1a nueva(o)", fecha 2023-10-19 (serial 45218)
$row95 = @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 45218, 15, 100112045, "Zapallo", "Camote", "1a nueva(o)", 550, 730, 750, 737, "`$/kilo (volumen en unidades)", "Perú", 737, 1, "Hortaliza")
for ($i = 0; $i -lt $row95.Length; $i++) {
    $ws.Cells.Item(95, $i + 1).Value = $row95[$i]
}

# New row 96: Camote, "2a nueva(o)", fecha 2023-10-19 (serial 45218)
$row96 = @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 45218, 15, 100112045, "Zapallo", "Camote", "2a nueva(o)", 180, 700, 720, 711, "`$/kilo (volumen en unidades)", "Perú", 711, 1, "Hortaliza")
for ($i = 0; $i -lt $row96.Length; $i++) {
    $ws.Cells.Item(96, $i + 1).Value = $row96[$i]
}
